$d = $word.ActiveDocument

# Locate the "30.201-5 Waiver." heading paragraph and split it so that a new
# bookmark "P30_201_5" wraps just the "30.201-5 " portion (mirroring the
# hyperlink anchor used elsewhere in the document / TOC).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "30.201-5 Waiver.`r") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $pStart = $target.Range.Start
    $splitOffset = $pStart + 9   # length of "30.201-5 "
    $bmRange = $d.Range($pStart, $splitOffset)
    $d.Bookmarks.Add("P30_201_5", $bmRange)
}
